# Generate Report for Handoff
#
# - Status moves from "In Translation" to "Ready for handoff" (Overview!E2,
#   Overview!F2, zh-cn!C2, de-de!C2 all share this text).
# - The handoff timestamps advance a few seconds:
#     Overview!G2 and de-de!H2 share "2016-08-19 02:37:04"
#     zh-cn!H2 becomes "2016-08-19 02:36:57"
# - The "Status" columns widen to fit the new, longer text.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

# --- Status text: "In Translation" -> "Ready for handoff" ---
$overview.Range("E2").Value = "Ready for handoff"
$overview.Range("F2").Value = "Ready for handoff"
$zhcn.Range("C2").Value = "Ready for handoff"
$dede.Range("C2").Value = "Ready for handoff"

# --- Handoff datetimes ---
$overview.Range("G2").Value = "2016-08-19 02:37:04"
$zhcn.Range("H2").Value = "2016-08-19 02:36:57"
$dede.Range("H2").Value = "2016-08-19 02:37:04"

# --- Widen the Status columns to fit "Ready for handoff" ---
$overview.Columns.Item(5).ColumnWidth = 16.33
$overview.Columns.Item(6).ColumnWidth = 16.33
$zhcn.Columns.Item(3).ColumnWidth = 16.33
$dede.Columns.Item(3).ColumnWidth = 16.33
